{"js": "const replacements = [\n  [\"26\u00f78=\", \"65\u00f79=\"],\n  [\"52\u00f73=\", \"59\u00f78=\"],\n  [\"12\u00f76=\", \"96\u00f77=\"],\n  [\"20\u00f79=\", \"42\u00f77=\"],\n  [\"51\u00f78=\", \"20\u00f78=\"],\n  [\"12\u00f79=\", \"69\u00f76=\"],\n  [\"24\u00f75=\", \"31\u00f73=\"],\n  [\"67\u00f79=\", \"94\u00f73=\"],\n  [\"98\u00f79=\", \"64\u00f74=\"],\n  [\"34\u00f75=\", \"74\u00f76=\"],\n  [\"56\u00f77=\", \"15\u00f78=\"],\n  [\"58\u00f78=\", \"66\u00f74=\"],\n  [\"65\u00f73=\", \"54\u00f73=\"],\n  [\"23\u00f77=\", \"89\u00f76=\"],\n  [\"39\u00f74=\", \"36\u00f75=\"],\n  [\"33\u00f75=\", \"48\u00f76=\"],\n  [\"93\u00f73=\", \"17\u00f74=\"],\n  [\"19\u00f75=\", \"90\u00f78=\"],\n  [\"59\u00f72=\", \"42\u00f79=\"],\n  [\"29\u00f76=\", \"41\u00f78=\"],\n  [\"32\u00f74=\", \"97\u00f74=\"],\n  [\"42\u00f74=\", \"89\u00f72=\"],\n  [\"84\u00f78=\", \"92\u00f74=\"],\n  [\"93\u00f74=\", \"78\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"26\u00f78=\", \"65\u00f79=\"),\n    @(\"52\u00f73=\", \"59\u00f78=\"),\n    @(\"12\u00f76=\", \"96\u00f77=\"),\n    @(\"20\u00f79=\", \"42\u00f77=\"),\n    @(\"51\u00f78=\", \"20\u00f78=\"),\n    @(\"12\u00f79=\", \"69\u00f76=\"),\n    @(\"24\u00f75=\", \"31\u00f73=\"),\n    @(\"67\u00f79=\", \"94\u00f73=\"),\n    @(\"98\u00f79=\", \"64\u00f74=\"),\n    @(\"34\u00f75=\", \"74\u00f76=\"),\n    @(\"56\u00f77=\", \"15\u00f78=\"),\n    @(\"58\u00f78=\", \"66\u00f74=\"),\n    @(\"65\u00f73=\", \"54\u00f73=\"),\n    @(\"23\u00f77=\", \"89\u00f76=\"),\n    @(\"39\u00f74=\", \"36\u00f75=\"),\n    @(\"33\u00f75=\", \"48\u00f76=\"),\n    @(\"93\u00f73=\", \"17\u00f74=\"),\n    @(\"19\u00f75=\", \"90\u00f78=\"),\n    @(\"59\u00f72=\", \"42\u00f79=\"),\n    @(\"29\u00f76=\", \"41\u00f78=\"),\n    @(\"32\u00f74=\", \"97\u00f74=\"),\n    @(\"42\u00f74=\", \"89\u00f72=\"),\n    @(\"84\u00f78=\", \"92\u00f74=\"),\n    @(\"93\u00f74=\", \"78\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
